# Updated symbol list on Thu Dec 29 16:50:07 UTC 2022 with GitHub Actions
#
# This refreshes the crypto price/volume feed on Sheet1: most rows just get
# an updated Price (column D), rows 18-24 are a coin-ranking reshuffle
# (each coin's Coin/Link/Price/Volume moved up one row, with "One" wrapping
# around to the bottom), and a handful of Volume(1h) labels (column E) drop
# or pick up a "Bestin24h"/"Worstin24h" suffix.
#
# All the cells touched here are stored as plain text in the workbook (not
# numbers), so for any cell whose new value happens to look numeric we force
# the Text format first -- otherwise Excel would silently convert it to a
# number and normalise away the exact decimal formatting (trailing zeros,
# leading zeros, etc.) that the source data relies on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.40"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.289"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05820"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.474"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.129"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8169"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8791"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1379"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06991"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03121"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02944"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09403"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.746"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001524"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04702"

# Row 18: was "One" -> now "TigerCash" (shifted up from row 19)
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006242"
$ws.Range("E18").Value = "17TigerCashTCH"

# Row 19: was "TigerCash" -> now "BitKan" (shifted up from row 20)
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001239"
$ws.Range("E19").Value = "18BitKanKAN"

# Row 20: was "BitKan" -> now "HotbitToken" (shifted up from row 21)
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004674"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21: was "HotbitToken" -> now "NitroEx" (shifted up from row 22)
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00006100"
$ws.Range("E21").Value = "20NitroExNTXWorstin24h"

# Row 22: was "NitroEx" -> now "LEO" (shifted up from row 23)
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.534"
$ws.Range("E22").Value = "21LEOLEO"

# Row 23: was "LEO" -> now "BTSEToken" (shifted up from row 24)
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.152"
$ws.Range("E23").Value = "22BTSETokenBTSE"

# Row 24: was "BTSEToken" -> now "One" (wrapped around from row 18)
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01006"
$ws.Range("E24").Value = "23OneONEBestin24h"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3182"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002332"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03725"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006470"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1058"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003400"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008272"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005256"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3695"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002415"
$ws.Range("E48").Value = "47BOLOBOLO"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
